# Petty cash book - end of day update, 27/28-Jun-2021
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the day's transaction detail rows (dates 21-Jun through 25-Jun-2021),
# keeping the running-balance formulas in column E intact.
$ws.Range("A4:D33").ClearContents()
$ws.Range("D3").ClearContents()

# Roll A3 forward to the new day being recorded.
$ws.Range("A3").Value = 44375

# Update the opening balance for the day.
$ws.Range("E2").Value = 479225

# Move the frozen-pane view back to the top of the data and park the
# selection on the next entry cell.
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A3").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("C4").Select()
